$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.480.34"
$ws.Range("E2").Value = "  -2.77%  "
$ws.Range("D3").Value = "2.466.69"
$ws.Range("E3").Value = "  -2.37%  "
$ws.Range("E4").Value = "  +0.96%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.26"
$ws.Range("E5").Value = "  -0.69%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "91.07"
$ws.Range("E6").Value = "  -7.55%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.540"
$ws.Range("E7").Value = "  -3.91%  "
$ws.Range("E8").Value = "  +0.86%  "
$ws.Range("E9").Value = "  -5.12%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "32.63"
$ws.Range("E10").Value = "  -7.16%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0775"
$ws.Range("E11").Value = "  -3.04%  "
$ws.Range("E12").Value = "  -0.13%  "
$ws.Range("D13").Value = "2.843.44"
$ws.Range("E13").Value = "  -2.47%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.81"
$ws.Range("E14").Value = "  -5.54%  "
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.17"
$ws.Range("E15").Value = "  -0.19%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "2.422.29"
$ws.Range("E16").Value = "  -4.43%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.777"
$ws.Range("E17").Value = "  -3.88%  "
$ws.Range("D18").Value = "41.233.08"
$ws.Range("E18").Value = "  -3.32%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.25"
$ws.Range("E19").Value = "  -4.96%  "
$ws.Range("D20").Value = "0.0₃0915"
$ws.Range("E20").Value = "  -2.51%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "70.28"
$ws.Range("E21").Value = "  +1.86%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.91"
$ws.Range("E22").Value = "  -9.97%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "233.92"
$ws.Range("E23").Value = "  -3.14%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.71"
$ws.Range("E25").Value = "  +0.18%  "
$ws.Range("E26").Value = "  -6.37%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "23.82"
$ws.Range("E27").Value = "  -6.60%  "
$ws.Range("E28").Value = "  -0.79%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.65"
$ws.Range("E29").Value = "  -3.39%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "35.88"
$ws.Range("E30").Value = "  -4.33%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "152.39"
$ws.Range("E31").Value = "  -2.02%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.39"
$ws.Range("E32").Value = "  -8.47%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.57"
$ws.Range("E33").Value = "  -5.50%  "
$ws.Range("E34").Value = "  -3.60%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0745"
$ws.Range("E35").Value = "  -4.63%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.34"
$ws.Range("E36").Value = "  -1.41%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.97"
$ws.Range("E37").Value = "  -4.86%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.83"
$ws.Range("E38").Value = "  -7.02%  "
$ws.Range("E39").Value = "  -3.80%  "
$ws.Range("E40").Value = "  -8.23%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.01"
$ws.Range("E41").Value = "  -4.80%  "
$ws.Range("E42").Value = "  +1.20%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "19.16"
$ws.Range("E43").Value = "  -12.10%  "
$ws.Range("D44").Value = "1.954.81"
$ws.Range("E44").Value = "  -3.60%  "
$ws.Range("E45").Value = "  -4.99%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.92"
$ws.Range("E46").Value = "  -9.01%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.59"
$ws.Range("E47").Value = "  -2.77%  "
$ws.Range("D48").Value = "2.709.95"
$ws.Range("E48").Value = "  -2.20%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "95.03"
$ws.Range("E49").Value = "  -5.07%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "67.28"
$ws.Range("E50").Value = "  -6.37%  "
$ws.Range("E51").Value = "  -7.08%  "
